$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row (row 1) from the placeholder "q"/"a" labels
# to the full "question"/"answer" labels used for rendering survey context.
$ws.Range("A1").Value = "question"
$ws.Range("B1").Value = "answer"

# Move the active selection to A2, matching the refreshed view state.
$ws.Range("A2").Select()
